# Update the dSF (column F) values for specific rows as per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -8
    3  = -3
    4  = 2
    8  = -3
    10 = -8
    12 = -5
    13 = -10
    16 = 2
    17 = -8
    18 = -6
    19 = -5
    21 = -7
    22 = -4
    25 = 3
    27 = -2
    28 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
